# Insert a new data row at row 44, shifting existing rows 44:143 down to 45:144,
# then populate the new row 44 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 44 and below down by one row.
$ws.Rows.Item(44).Insert()

# Populate the newly inserted row 44 with the new record.
$ws.Cells.Item(44, 1).Value  = 1
$ws.Cells.Item(44, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(44, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(44, 4).Value  = 45076
$ws.Cells.Item(44, 5).Value  = 15
$ws.Cells.Item(44, 6).Value  = "Fruta"
$ws.Cells.Item(44, 7).Value  = 100102
$ws.Cells.Item(44, 8).Value  = "Cítricos"
$ws.Cells.Item(44, 9).Value  = 100102005
$ws.Cells.Item(44, 10).Value = "Naranja"
$ws.Cells.Item(44, 11).Value = "Fukumoto"
$ws.Cells.Item(44, 12).Value = "Primera"
$ws.Cells.Item(44, 13).Value = 350
$ws.Cells.Item(44, 14).Value = 1200
$ws.Cells.Item(44, 15).Value = 1250
$ws.Cells.Item(44, 16).Value = 1221
$ws.Cells.Item(44, 17).Value = "$/kilo (en caja de 20 kilos)"
$ws.Cells.Item(44, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(44, 19).Value = 1221
$ws.Cells.Item(44, 20).Value = 1
